# Generate Report for Handoff
#
# Moves the localization-status report from "In Translation" to
# "Ready for handoff" and refreshes the associated timestamps on all
# three sheets (Overview, zh-cn, de-de). Also widens the status columns
# that now hold the longer "Ready for handoff" label.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-06 10:52:21"

# Status text grew ("In Translation" -> "Ready for handoff"); widen the
# zh-cn / de-de status columns to fit the new text.
$overview.Range("E1").ColumnWidth = 16.29
$overview.Range("F1").ColumnWidth = 16.29

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-06 10:52:16"
$zhcn.Range("C1").ColumnWidth = 16.29

# --- de-de sheet ----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-06 10:52:21"
$dede.Range("C1").ColumnWidth = 16.29
